$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows
# (2 through 151) from 45192 (2023-09-23) to 45202 (2023-10-03).
for ($r = 2; $r -le 151; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}
